# Denmark Division 2 - base update (29-03-2024 17:05)
# The underlying scraped match rows were re-keyed; for a handful of rows the
# "id"/metadata column (A) stayed put on its row, but all the match data
# (columns B..AC: match id, teams, scores, odds, etc.) needed to be rotated
# among the affected rows to line up with the corrected source order.
#
# For each group below, the data in columns B..AC is cyclically shifted up
# by one row (row[i] <- row[i+1], with the last row in the group wrapping
# around and receiving the first row's original data). Column A (the running
# index) is left untouched on every row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 2   # column B
$lastCol  = 29  # column AC

$rowGroups = @(
    @(13, 14, 15),
    @(19, 20),
    @(29, 30),
    @(82, 83, 85),
    @(113, 114),
    @(118, 119),
    @(127, 128)
)

foreach ($group in $rowGroups) {
    $n = $group.Length

    # Snapshot the current B..AC values for every row in this group first,
    # so the subsequent writes don't clobber data we still need to read.
    $snapshots = @()
    foreach ($r in $group) {
        $rowVals = @{}
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $rowVals[$c] = $ws.Cells.Item($r, $c).Value2
        }
        $snapshots += ,$rowVals
    }

    # Write back: row i gets the snapshot originally belonging to row i+1
    # (wrapping to row 0 for the last row in the group).
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $group[$i]
        $srcVals = $snapshots[($i + 1) % $n]
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $ws.Cells.Item($destRow, $c).Value = $srcVals[$c]
        }
    }
}
